$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.726.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.542.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.78"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.542.75"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.999.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.562.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.548.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.88"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.57"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0830"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.95"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +13.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "417.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.407"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.11"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.42"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.61%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.93"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.610"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0968"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.84"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.37%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.60%  "
